$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these price cells keep their original text formatting
# (Excel would otherwise auto-convert plain numeric-looking strings to numbers)
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '61.194.83'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '2.404.49'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  +0.49%  '
$ws.Range('D5').Value = '568.39'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').Value = '142.11'
$ws.Range('E6').Value = '  +1.22%  '
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('D8').Value = '0.536'
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('D9').Value = '2.413.45'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  +1.22%  '
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('E12').Value = '  +2.70%  '
$ws.Range('D13').Value = '0.346'
$ws.Range('E13').Value = '  +2.11%  '
$ws.Range('D14').Value = '26.44'
$ws.Range('E14').Value = '  +1.09%  '
$ws.Range('E15').Value = '  -0.63%  '
$ws.Range('D16').Value = '2.799.02'
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('D17').Value = '60.854.93'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = '2.409.64'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '8.07'
$ws.Range('E19').Value = '  +3.94%  '
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('D21').Value = '323.85'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('E23').Value = '  -0.70%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('B25').Value = 'SuiNetwork'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D25').Value = '1.93'
$ws.Range('E25').Value = '  +4.24%  '
$ws.Range('D26').Value = '65.20'
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('D27').Value = '595.83'
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('D28').Value = '8.27'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').Value = '0.0₃0948'
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('D30').Value = '2.519.16'
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('D31').Value = '7.99'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('E32').Value = '  +1.63%  '
$ws.Range('D33').Value = '1.80'
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.46'
$ws.Range('E36').Value = '  +2.71%  '
$ws.Range('E37').Value = '  +1.09%  '
$ws.Range('D38').Value = '4.62'
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('E39').Value = '  -0.63%  '
$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D40').Value = '18.34'
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').Value = '5.28'
$ws.Range('E41').Value = '  +1.84%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '1.70'
$ws.Range('E43').Value = '  +1.10%  '
$ws.Range('D44').Value = '2.52'
$ws.Range('E44').Value = '  +4.19%  '
$ws.Range('D45').Value = '42.02'
$ws.Range('E45').Value = '  +1.74%  '
$ws.Range('D46').Value = '0.0₆0284'
$ws.Range('E46').Value = '  -5.85%  '
$ws.Range('D47').Value = '141.48'
$ws.Range('E47').Value = '  -1.21%  '
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D49').Value = '19.88'
$ws.Range('E49').Value = '  +1.33%  '
$ws.Range('E50').Value = '  +0.52%  '
$ws.Range('D51').Value = '0.0509'
$ws.Range('E51').Value = '  +1.18%  '

# Restore default (no explicit number format) style on those cells
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
